# Weekly update: insert a new week's price record for Orégano (Vega Central
# Mapocho de Santiago) at the top of the data block (row 33), pushing the
# existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 33 (shifts rows 33:44 down to 34:45, carrying
# formatting/styles along for free, same as Excel's own Insert command).
$ws.Rows(33).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A33").Value2 = 9
$ws.Range("B33").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value2 = "Metropolitana"
$ws.Range("D33").Value2 = 44559
$ws.Range("E33").Value2 = 13
$ws.Range("F33").Value2 = 100112029
$ws.Range("G33").Value2 = "Orégano"
$ws.Range("H33").Value2 = "Sin especificar"
$ws.Range("I33").Value2 = "Primera"
$ws.Range("J33").Value2 = 7
$ws.Range("K33").Value2 = 10000
$ws.Range("L33").Value2 = 12000
$ws.Range("M33").Value2 = 11143
$ws.Range("N33").Value2 = "`$/docena de atados"
$ws.Range("O33").Value2 = "Región Metropolitana"
$ws.Range("P33").Value2 = 3714
$ws.Range("Q33").Value2 = 3
$ws.Range("R33").Value2 = "Hortaliza"
